# Commit: Update works images 2026-01-21 17:22:18
#
# A name ("zd9ff0d0") that was sitting as the next-available entry in the
# "Sheet1" pool of unused random names gets consumed: it is removed from
# the pool (row 2 deleted, all later rows shift up one, shrinking the
# range from A1:A460 to A1:A459) and recorded as newly "used" on the
# "used" sheet, where a row is appended with that id, the source filename
# and the timestamp it was used at (growing that sheet from A1:C39 to
# A1:C40).

$wb = $excel.ActiveWorkbook

$namesSheet = $wb.Worksheets.Item("Sheet1")
$usedSheet  = $wb.Worksheets.Item("used")

# The id about to be used is the one currently sitting in the pool's
# second row (row 1 is already consumed/ahead in the queue).
$consumedId = $namesSheet.Cells.Item(2, 1).Value2

# Remove it from the pool - deleting the row shifts everything below up,
# shrinking the sheet from 460 to 459 rows.
$namesSheet.Rows.Item(2).Delete()

# Append the new "used" record right after the last existing row (39).
$nextRow = $usedSheet.UsedRange.Rows.Count + 1

$usedSheet.Cells.Item($nextRow, 1).Value = $consumedId
$usedSheet.Cells.Item($nextRow, 2).Value = "ChatGPT Image 2026年1月21日 17_21_37.png"
$usedSheet.Cells.Item($nextRow, 3).Value = "2026-01-21 17:22:11"
